# Fix pluricandidature nello stesso collegio
# Updates polling data on "liste_naz" sheet (new "Insieme per il futuro" party
# replacing "Noi con l'Italia", refreshed percentages, reordered rows, a new
# "Altri 6" row) and the source/value on "altri_dati" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # liste_naz
$ws2 = $wb.Worksheets.Item(2)   # altri_dati

# ---------------------------------------------------------------------------
# liste_naz: rows 2 and 7 are unchanged. Rows 3,4,5,6,8 only get new
# PERCENTUALE values (columns C and F). Rows 9-12 are fully rewritten because
# the list got re-sorted by percentage and a party swap happened. Rows 13-16
# are unchanged. Row 17 gets a new percentage and row 18 is a brand-new row.
# ---------------------------------------------------------------------------

# Row 3: Partito Democratico 22.1% -> 22.5%
$ws1.Range("C3").Value = 0.225
$ws1.Range("F3").Value = 0.225

# Row 4: Lega 14.0% -> 13.4%
$ws1.Range("C4").Value = 0.134
$ws1.Range("F4").Value = 0.134

# Row 5: Movimento 5 Stelle 11.2% -> 9.8%
$ws1.Range("C5").Value = 0.098
$ws1.Range("F5").Value = 0.098

# Row 6: Forza Italia 7.4% -> 8.3%
$ws1.Range("C6").Value = 0.083
$ws1.Range("F6").Value = 0.083

# Row 8: Europa Verde - Sinistra Italiana 3.8% -> 4.2%
$ws1.Range("C8").Value = 0.042
$ws1.Range("F8").Value = 0.042

# Row 9: now "Insieme per il futuro" (IxF), replacing what used to be here
$ws1.Range("A9").Value = "Insieme per il futuro"
$ws1.Range("B9").Value = "IxF"
$ws1.Range("C9").Value = 0.026
$ws1.Range("D9").Value = "SX"
$ws1.Range("F9").Value = 0.026
$ws1.Range("G9").Value = 0.23
$ws1.Range("H9").Value = $true
$ws1.Range("I9").Value = 55
$ws1.Range("J9").Value = "M5S"

# Row 10: now Italia Viva
$ws1.Range("A10").Value = "Italia Viva"
$ws1.Range("B10").Value = "IV"
$ws1.Range("C10").Value = 0.018
$ws1.Range("D10").Value = "SX"
$ws1.Range("F10").Value = 0.018
$ws1.Range("G10").Value = 0.35
$ws1.Range("H10").Value = $false
$ws1.Range("I10").Value = 300
$ws1.Range("J10").Value = "CENTRO"

# Row 11: now Italexit (and loses its MINORANZA/D value)
$ws1.Range("A11").Value = "Italexit"
$ws1.Range("B11").Value = "Italexit"
$ws1.Range("C11").Value = 0.02
$ws1.Range("D11").ClearContents()
$ws1.Range("F11").Value = 0.02
$ws1.Range("G11").Value = 0.35
$ws1.Range("H11").Value = $false
$ws1.Range("I11").Value = 270
$ws1.Range("J11").Value = "DX"

# Row 12: now Art. 1 - MDP (replaces the removed "Noi con l'Italia" row),
# gains an ABBREV (MDP) and a K style cell
$ws1.Range("A12").Value = "Art. 1 - MDP"
$ws1.Range("B12").Value = "MDP"
$ws1.Range("C12").Value = 0.016
$ws1.Range("D12").Value = "SX"
$ws1.Range("F12").Value = 0.016
$ws1.Range("G12").Value = 0.4
$ws1.Range("H12").Value = $false
$ws1.Range("I12").Value = 15
$ws1.Range("J12").Value = "CENTRO"
$ws1.Range("K12").NumberFormat = $ws1.Range("K11").NumberFormat()

# Row 17: Altri 5 now 1.0% (was 0.3%)
$ws1.Range("C17").Value = 0.01
$ws1.Range("F17").Value = 0.01

# Row 18: brand new "Altri 6" row
$ws1.Range("A18").Value = "Altri 6"
$ws1.Range("C18").Value = 0.001
$ws1.Range("E18").Value = $false
$ws1.Range("F18").Value = 0.001
$ws1.Range("G18").Value = 0.25
$ws1.Range("H18").Value = $false
$ws1.Range("J18").Value = "ALTRI"
$ws1.Range("C18:F18").NumberFormat = $ws1.Range("C17:F17").NumberFormat()
$ws1.Range("G18").NumberFormat = $ws1.Range("G17").NumberFormat()

# ---------------------------------------------------------------------------
# altri_dati: update the poll source label and the abstention value
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Youtrend per Sky TG24 25/7"
$ws2.Range("B2").Value = 0.417
$ws2.Columns.Item(1).AutoFit()

# liste_naz is the sheet the author left active in this revision
$ws1.Activate()
